$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values in column F (dSF) reflecting repulled/pushed data and mean recalculation
$updates = @{
    4  = -2
    7  = 7
    9  = 1
    16 = -2
    20 = -2
    22 = -1
    25 = 6
    26 = -1
    27 = 1
    35 = 6
    36 = 1
    37 = -4
    41 = -2
    42 = -1
    44 = -1
    47 = 2
    49 = 0
    51 = -2
    52 = -3
    58 = 0
    59 = -2
    60 = 3
    62 = 4
    63 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
